$p = $ppt.ActivePresentation

# Slide 10 (sldId 274) holds the "Results and Analysis" bullet list that
# incorrectly called out "linear regression" instead of "logistic regression".
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item("Content Placeholder 9")

$tf = $shp.TextFrame
$tr = $tf.TextRange

# The first paragraph ("The two linear regression models each exceeded 95%
# accuracy.") is a single run - update just that run's text in place so the
# rest of the text frame (other paragraphs/runs) is left untouched.
$run = $tr.Runs(1)
$run.Text = "The two logistic regression models each exceeded 95% accuracy."
